# Sync attendance_reports: fix ordering of the "Recorded By" names in column G.
#
# The report exporter populates the "Recorded By" cell from an unordered
# collection of recorder identities, so re-running the export can change the
# order the names are listed in. This resyncs column G on the "Session
# Analysis Results" sheet so the affected rows show the recorder names in
# the canonical order (first two names swapped; any trailing name such as
# "System" in a 3-name list stays put).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$colRecordedBy = 7   # column G ("Recorded By")

# Exact "before" -> "after" text swaps that need to be applied, scoped to
# only the cells that actually hold one of these values (other
# two/three-name combinations, e.g. "backup@backdoor.com, System" or
# "admin@admin.com, System", are left completely untouched).
$renameMap = @{
    "system, backup@backdoor.com, System" = "backup@backdoor.com, system, System"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
}

$changed = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colRecordedBy)
    $text = $cell.Value2

    if ($null -eq $text) { continue }

    if ($renameMap.ContainsKey($text)) {
        $newText = $renameMap[$text]
        if ($newText -ne $text) {
            $cell.Value = $newText
            $changed++
        }
    }
}

Write-Host "Resynced 'Recorded By' ordering on $changed rows."
